$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: add Kit Faltando value
$ws.Range("I6").Value = "2 módulo internet"

# Row 7: fill in ID, Cliente, Descricao, Status (Técnico and Kit Faltando already set)
$ws.Range("B7").Value = "'0554"
$ws.Range("C7").Value = "Valinhos Departamento"
$ws.Range("D7").Value = "Moisés pedindo ajuste de sensibilidade em um sensor."
$ws.Range("G7").Value = "Pendente"

# Row 8: new entry
$ws.Range("A8").Value = "Roberto"
$ws.Range("B8").Value = "'0706"
$ws.Range("C8").Value = "Lar das Meninas"
$ws.Range("D8").Value = "Várias câmeras sem imagem."
$ws.Range("G8").Value = "Pendente"

# Row 9: new entry
$ws.Range("A9").Value = "Roberto"
$ws.Range("B9").Value = "'0756"
$ws.Range("C9").Value = "Caixa Escolar Manoel Correia"
$ws.Range("D9").Value = "Sem comunicação geral."
$ws.Range("G9").Value = "Pendente"

# Update sheet view: scroll and selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("H9").Select()
